# Auto-generated data update for violent-crime-full-year workbook
# Commit: Add data for 2025-07-16
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "Citywide Totals"; Cell = "L2"; Value = 3573 },
    @{ Sheet = "Citywide Totals"; Cell = "L3"; Value = 3742 },
    @{ Sheet = "Citywide Totals"; Cell = "E4"; Value = 2052 },
    @{ Sheet = "Citywide Totals"; Cell = "L4"; Value = 929 },
    @{ Sheet = "Citywide Totals"; Cell = "L5"; Value = 222 },
    @{ Sheet = "Citywide Totals"; Cell = "L6"; Value = 3263 },
    @{ Sheet = "Citywide Totals"; Cell = "E7"; Value = 26057 },
    @{ Sheet = "Citywide Totals"; Cell = "L7"; Value = 11729 },
    @{ Sheet = "Logan Square"; Cell = "L2"; Value = 40 },
    @{ Sheet = "Logan Square"; Cell = "L4"; Value = 12 },
    @{ Sheet = "Logan Square"; Cell = "L7"; Value = 130 },
    @{ Sheet = "Austin"; Cell = "L2"; Value = 218 },
    @{ Sheet = "Austin"; Cell = "L3"; Value = 253 },
    @{ Sheet = "Austin"; Cell = "L4"; Value = 50 },
    @{ Sheet = "Austin"; Cell = "L7"; Value = 758 },
    @{ Sheet = "Garfield Park"; Cell = "L2"; Value = 153 },
    @{ Sheet = "Garfield Park"; Cell = "L3"; Value = 181 },
    @{ Sheet = "Garfield Park"; Cell = "L6"; Value = 181 },
    @{ Sheet = "Garfield Park"; Cell = "L7"; Value = 553 },
    @{ Sheet = "West Pullman"; Cell = "L2"; Value = 65 },
    @{ Sheet = "West Pullman"; Cell = "E4"; Value = 24 },
    @{ Sheet = "West Pullman"; Cell = "L4"; Value = 9 },
    @{ Sheet = "West Pullman"; Cell = "E7"; Value = 441 },
    @{ Sheet = "New City"; Cell = "L3"; Value = 66 },
    @{ Sheet = "New City"; Cell = "L6"; Value = 64 },
    @{ Sheet = "Woodlawn"; Cell = "L5"; Value = 5 },
    @{ Sheet = "Woodlawn"; Cell = "L7"; Value = 198 },
    @{ Sheet = "Fuller Park"; Cell = "L6"; Value = 26 },
    @{ Sheet = "Fuller Park"; Cell = "L7"; Value = 59 },
    @{ Sheet = "By Neighborhood"; Cell = "L6"; Value = 92 },
    @{ Sheet = "By Neighborhood"; Cell = "L8"; Value = 758 },
    @{ Sheet = "By Neighborhood"; Cell = "L11"; Value = 195 },
    @{ Sheet = "By Neighborhood"; Cell = "L19"; Value = 335 },
    @{ Sheet = "By Neighborhood"; Cell = "L27"; Value = 107 },
    @{ Sheet = "By Neighborhood"; Cell = "L29"; Value = 632 },
    @{ Sheet = "By Neighborhood"; Cell = "L30"; Value = 59 },
    @{ Sheet = "By Neighborhood"; Cell = "L31"; Value = 113 },
    @{ Sheet = "By Neighborhood"; Cell = "L33"; Value = 553 },
    @{ Sheet = "By Neighborhood"; Cell = "L36"; Value = 159 },
    @{ Sheet = "By Neighborhood"; Cell = "L40"; Value = 35 },
    @{ Sheet = "By Neighborhood"; Cell = "L41"; Value = 52 },
    @{ Sheet = "By Neighborhood"; Cell = "L42"; Value = 374 },
    @{ Sheet = "By Neighborhood"; Cell = "L43"; Value = 87 },
    @{ Sheet = "By Neighborhood"; Cell = "L50"; Value = 57 },
    @{ Sheet = "By Neighborhood"; Cell = "L53"; Value = 130 },
    @{ Sheet = "By Neighborhood"; Cell = "L55"; Value = 111 },
    @{ Sheet = "By Neighborhood"; Cell = "L63"; Value = 41 },
    @{ Sheet = "By Neighborhood"; Cell = "L67"; Value = 418 },
    @{ Sheet = "By Neighborhood"; Cell = "L73"; Value = 99 },
    @{ Sheet = "By Neighborhood"; Cell = "L76"; Value = 173 },
    @{ Sheet = "By Neighborhood"; Cell = "L77"; Value = 72 },
    @{ Sheet = "By Neighborhood"; Cell = "L78"; Value = 146 },
    @{ Sheet = "By Neighborhood"; Cell = "L79"; Value = 306 },
    @{ Sheet = "By Neighborhood"; Cell = "L84"; Value = 117 },
    @{ Sheet = "By Neighborhood"; Cell = "L85"; Value = 602 },
    @{ Sheet = "By Neighborhood"; Cell = "L88"; Value = 132 },
    @{ Sheet = "By Neighborhood"; Cell = "L89"; Value = 165 },
    @{ Sheet = "By Neighborhood"; Cell = "L90"; Value = 112 },
    @{ Sheet = "By Neighborhood"; Cell = "L91"; Value = 167 },
    @{ Sheet = "By Neighborhood"; Cell = "L93"; Value = 63 },
    @{ Sheet = "By Neighborhood"; Cell = "L94"; Value = 138 },
    @{ Sheet = "By Neighborhood"; Cell = "E95"; Value = 441 },
    @{ Sheet = "By Neighborhood"; Cell = "L99"; Value = 198 },
    @{ Sheet = "By Neighborhood"; Cell = "E101"; Value = 26057 },
    @{ Sheet = "By Neighborhood"; Cell = "L101"; Value = 11729 },
    @{ Sheet = "Gage Park"; Cell = "L3"; Value = 32 },
    @{ Sheet = "Gage Park"; Cell = "L7"; Value = 113 },
    @{ Sheet = "North Lawndale"; Cell = "L4"; Value = 31 },
    @{ Sheet = "North Lawndale"; Cell = "L6"; Value = 98 },
    @{ Sheet = "North Lawndale"; Cell = "L7"; Value = 418 },
    @{ Sheet = "South Deering"; Cell = "L2"; Value = 40 },
    @{ Sheet = "South Deering"; Cell = "L6"; Value = 27 },
    @{ Sheet = "South Deering"; Cell = "L7"; Value = 117 },
    @{ Sheet = "Englewood"; Cell = "L2"; Value = 189 },
    @{ Sheet = "Englewood"; Cell = "L3"; Value = 243 },
    @{ Sheet = "Englewood"; Cell = "L6"; Value = 158 },
    @{ Sheet = "Englewood"; Cell = "L7"; Value = 632 },
    @{ Sheet = "Chatham"; Cell = "L2"; Value = 118 },
    @{ Sheet = "Chatham"; Cell = "L4"; Value = 15 },
    @{ Sheet = "Chatham"; Cell = "L7"; Value = 335 },
    @{ Sheet = "River North"; Cell = "L6"; Value = 80 },
    @{ Sheet = "River North"; Cell = "L7"; Value = 173 },
    @{ Sheet = "Ashburn"; Cell = "L6"; Value = 20 },
    @{ Sheet = "Ashburn"; Cell = "L7"; Value = 92 },
    @{ Sheet = "Hermosa"; Cell = "L6"; Value = 13 },
    @{ Sheet = "Hermosa"; Cell = "L7"; Value = 52 },
    @{ Sheet = "Humboldt Park"; Cell = "L6"; Value = 107 },
    @{ Sheet = "Humboldt Park"; Cell = "L7"; Value = 374 },
    @{ Sheet = "Rogers Park"; Cell = "L3"; Value = 41 },
    @{ Sheet = "Rogers Park"; Cell = "L7"; Value = 146 },
    @{ Sheet = "Lower West Side"; Cell = "L6"; Value = 29 },
    @{ Sheet = "Lower West Side"; Cell = "L7"; Value = 111 },
    @{ Sheet = "Washington Park"; Cell = "L2"; Value = 59 },
    @{ Sheet = "Washington Park"; Cell = "L7"; Value = 167 },
    @{ Sheet = "Roseland"; Cell = "L3"; Value = 111 },
    @{ Sheet = "Roseland"; Cell = "L6"; Value = 60 },
    @{ Sheet = "Roseland"; Cell = "L7"; Value = 306 },
    @{ Sheet = "Grand Boulevard"; Cell = "L2"; Value = 59 },
    @{ Sheet = "Grand Boulevard"; Cell = "L7"; Value = 159 },
    @{ Sheet = "West Lawn"; Cell = "L3"; Value = 17 },
    @{ Sheet = "West Lawn"; Cell = "L7"; Value = 63 },
    @{ Sheet = "West Loop"; Cell = "L4"; Value = 21 },
    @{ Sheet = "West Loop"; Cell = "L7"; Value = 138 },
    @{ Sheet = "Lincoln Square"; Cell = "L2"; Value = 23 },
    @{ Sheet = "Lincoln Square"; Cell = "L7"; Value = 57 },
    @{ Sheet = "Belmont Cragin"; Cell = "L2"; Value = 75 },
    @{ Sheet = "Belmont Cragin"; Cell = "L7"; Value = 195 },
    @{ Sheet = "Portage Park"; Cell = "L3"; Value = 28 },
    @{ Sheet = "Portage Park"; Cell = "L7"; Value = 99 },
    @{ Sheet = "United Center"; Cell = "L3"; Value = 48 },
    @{ Sheet = "United Center"; Cell = "L7"; Value = 132 },
    @{ Sheet = "Galewood"; Cell = "L2"; Value = 8 },
    @{ Sheet = "Galewood"; Cell = "L3"; Value = 3 },
    @{ Sheet = "Uptown"; Cell = "L2"; Value = 48 },
    @{ Sheet = "Uptown"; Cell = "L3"; Value = 47 },
    @{ Sheet = "Uptown"; Cell = "L7"; Value = 165 },
    @{ Sheet = "Edgewater"; Cell = "L6"; Value = 31 },
    @{ Sheet = "Edgewater"; Cell = "L7"; Value = 107 },
    @{ Sheet = "Washington Heights"; Cell = "L6"; Value = 30 },
    @{ Sheet = "Washington Heights"; Cell = "L7"; Value = 112 },
    @{ Sheet = "Hyde Park"; Cell = "L4"; Value = 16 },
    @{ Sheet = "Hyde Park"; Cell = "L7"; Value = 87 },
    @{ Sheet = "South Shore"; Cell = "L2"; Value = 174 },
    @{ Sheet = "South Shore"; Cell = "L3"; Value = 250 },
    @{ Sheet = "South Shore"; Cell = "L6"; Value = 126 },
    @{ Sheet = "South Shore"; Cell = "L7"; Value = 602 },
    @{ Sheet = "Riverdale"; Cell = "L2"; Value = 22 },
    @{ Sheet = "Riverdale"; Cell = "L3"; Value = 30 },
    @{ Sheet = "Riverdale"; Cell = "L7"; Value = 72 },
    @{ Sheet = "Hegewisch"; Cell = "L3"; Value = 15 },
    @{ Sheet = "Hegewisch"; Cell = "L7"; Value = 35 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
